$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the date-like columns (H = Date of Birth, J = Registration Date) stay
# plain text rather than being auto-converted to date serials by Excel.
$ws.Range("H2:H4").NumberFormat = "@"
$ws.Range("J2:J4").NumberFormat = "@"

# Row 2 (Angela Matthews -> Lisa White)
$ws.Range("A2").Value = "Lisa White"
$ws.Range("B2").Value = "rachelcopeland@example.com"
$ws.Range("C2").Value = "+1-850-985-4880x5059"
$ws.Range("D2").Value = "580.766.6504x18211"
$ws.Range("E2").Value = "Sancheztown"
$ws.Range("F2").Value = 11
$ws.Range("H2").Value = "2013-11-05"
$ws.Range("I2").Value = "9651 Julie Glens, New Erin, WI 21351"
$ws.Range("J2").Value = "2024-01-24"
$ws.Range("K2").Value = "Maryland"

# Row 3 (Robert Smith -> Bryan Larson)
$ws.Range("A3").Value = "Bryan Larson"
$ws.Range("B3").Value = "waltermichael@example.org"
$ws.Range("C3").Value = "001-725-321-2467x0153"
$ws.Range("D3").Value = "(727)597-4050x23173"
$ws.Range("E3").Value = "Port William"
$ws.Range("F3").Value = 8
$ws.Range("G3").Value = "Male"
$ws.Range("H3").Value = "2007-07-17"
$ws.Range("I3").Value = "1043 Wise Trail, Hortonburgh, FM 55732"
$ws.Range("J3").Value = "2024-01-09"
$ws.Range("K3").Value = "Missouri"

# Row 4 (Elizabeth Craig -> Edward Cook)
$ws.Range("A4").Value = "Edward Cook"
$ws.Range("B4").Value = "matthew77@example.org"
$ws.Range("C4").Value = "001-880-510-6930"
$ws.Range("D4").Value = "339.479.2643x938"
$ws.Range("E4").Value = "Hillview"
$ws.Range("F4").Value = 2
$ws.Range("H4").Value = "2012-04-22"
$ws.Range("I4").Value = "Unit 2912 Box 9879, DPO AE 32503"
$ws.Range("J4").Value = "2024-03-29"
$ws.Range("K4").Value = "Georgia"
